# Applies the "Add files via upload" revision to the wheat-stress dataset:
#  - Relabel the two header cells (A1: "Cultivar " -> "Cultivar",
#    B1: "Stress" -> "Treatment")
#  - Normalize five cultivar names in column A from ALL CAPS to Title Case
#    (TBIO CALIBRE -> Tbio Calibre, TBIO DUQUE -> Tbio Duque,
#     TBIO SOSSEGO -> Tbio Sossego, ORS FEROZ -> ORS Feroz,
#     TBIO CONVICTO -> Tbio Convicto)
#  - Move the sheet's viewport/selection down to the last cultivar block
#    (A123:A133)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabeling ---------------------------------------------
$ws.Range("A1").Value = "Cultivar"
$ws.Range("B1").Value = "Treatment"

# --- Normalize ALL-CAPS cultivar names in column A ----------------------
$rename = @{
    "TBIO CALIBRE"  = "Tbio Calibre"
    "TBIO DUQUE"    = "Tbio Duque"
    "TBIO SOSSEGO"  = "Tbio Sossego"
    "ORS FEROZ"     = "ORS Feroz"
    "TBIO CONVICTO" = "Tbio Convicto"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 133 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($rename.ContainsKey($current)) {
        $cell.Value = $rename[$current]
    }
}

# --- Scroll/selection: land on the final cultivar's block ---------------
$win = $excel.ActiveWindow
$win.ScrollRow = 118
$win.ScrollColumn = 1
$ws.Range("A123:A133").Select()
